# Updated cryptos list on Mon Jan  1 03:32:18 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into the Price (D) column as literal text, even when
# it looks like a number (e.g. "0.999"), without leaving the cell's style
# pointing at a different (Text / quote-prefixed) cell format than the
# original "General" default. Excel's COM layer auto-coerces plain numeric-
# looking strings to real numbers and also stamps a quote-prefix style when
# a leading apostrophe is used to force text - resetting Style to "Normal"
# afterwards drops that extra style so the cell's style index is unchanged.
function Set-PriceText($addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = "'" + $text
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-PriceText "D2" "42.438.45"
$ws.Range("E2").Value = "  +0.87%  "

# Row 3 - Ethereum
Set-PriceText "D3" "2.282.15"
$ws.Range("E3").Value = "  +0.16%  "

# Row 4 - TetherUSD
Set-PriceText "D4" "0.999"
$ws.Range("E4").Value = "  -0.24%  "

# Row 5 - BNB
Set-PriceText "D5" "310.52"
$ws.Range("E5").Value = "  -2.52%  "

# Row 6 - Solana
Set-PriceText "D6" "103.44"
$ws.Range("E6").Value = "  +3.15%  "

# Row 7 - XRP
Set-PriceText "D7" "0.621"
$ws.Range("E7").Value = "  -0.73%  "

# Row 8 - USDC (D unchanged)
$ws.Range("E8").Value = "  -0.05%  "

# Row 9 - Cardano
Set-PriceText "D9" "0.597"
$ws.Range("E9").Value = "  -0.44%  "

# Row 10 - Avalanche
Set-PriceText "D10" "38.79"
$ws.Range("E10").Value = "  +0.11%  "

# Row 11 - Dogecoin
Set-PriceText "D11" "0.0896"
$ws.Range("E11").Value = "  -0.30%  "

# Row 12 - Polkadot
Set-PriceText "D12" "8.23"
$ws.Range("E12").Value = "  +0.42%  "

# Row 13 - TRON (D unchanged)
$ws.Range("E13").Value = "  +1.55%  "

# Row 14 - Polygon
Set-PriceText "D14" "0.974"
$ws.Range("E14").Value = "  +2.83%  "

# Row 15 - Chainlink
Set-PriceText "D15" "15.01"
$ws.Range("E15").Value = "  +0.14%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-PriceText "D16" "2.632.61"
$ws.Range("E16").Value = "  +0.21%  "

# Row 17 - WrappedEther
Set-PriceText "D17" "2.284.18"
$ws.Range("E17").Value = "  +0.11%  "

# Row 18 - WrappedBTC
Set-PriceText "D18" "42.381.06"
$ws.Range("E18").Value = "  +0.77%  "

# Row 19 - Uniswap
Set-PriceText "D19" "7.24"
$ws.Range("E19").Value = "  -0.70%  "

# Row 20 - ShibaInu : unchanged

# Row 21 - InternetComputer(DFINITY)
Set-PriceText "D21" "13.43"
$ws.Range("E21").Value = "  +6.12%  "

# Row 22 - Litecoin
Set-PriceText "D22" "72.82"
$ws.Range("E22").Value = "  +0.39%  "

# Row 23 - PancakeSwap
Set-PriceText "D23" "3.41"
$ws.Range("E23").Value = "  -3.28%  "

# Row 24 - BitcoinCash
Set-PriceText "D24" "261.85"
$ws.Range("E24").Value = "  -2.09%  "

# Row 25 - ImmutableX (D unchanged)
$ws.Range("E25").Value = "  -1.47%  "

# Row 26 - Dai (D unchanged)
$ws.Range("E26").Value = "  +0.35%  "

# Row 27 - Cosmos
Set-PriceText "D27" "10.65"
$ws.Range("E27").Value = "  -0.81%  "

# Row 28 - was Toncoin, now Filecoin (rows 28/29 swapped with new data)
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-PriceText "D28" "7.00"
$ws.Range("E28").Value = "  +16.55%  "

# Row 29 - was Filecoin, now Toncoin
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-PriceText "D29" "2.27"
$ws.Range("E29").Value = "  -2.21%  "

# Row 30 - EthereumClassic
Set-PriceText "D30" "22.19"
$ws.Range("E30").Value = "  -0.46%  "

# Row 31 - was Monero, now InjectiveProtocol (rows 31/32 swapped with new data)
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-PriceText "D31" "35.58"
$ws.Range("E31").Value = "  -4.13%  "

# Row 32 - was InjectiveProtocol, now Monero
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-PriceText "D32" "164.37"
$ws.Range("E32").Value = "  +0.34%  "

# Row 33 - Hedera
Set-PriceText "D33" "0.0856"
$ws.Range("E33").Value = "  -0.87%  "

# Row 34 - Stellar
Set-PriceText "D34" "0.129"

# Row 35 - WEMIXToken (D unchanged)
$ws.Range("E35").Value = "  +1.21%  "

# Row 36 - Kaspa (D unchanged)
$ws.Range("E36").Value = "  -2.39%  "

# Row 37 - RenderToken (D unchanged)
$ws.Range("E37").Value = "  -1.17%  "

# Row 38 - VeChain
Set-PriceText "D38" "0.0348"
$ws.Range("E38").Value = "  -1.14%  "

# Row 39 - NEARProtocol
Set-PriceText "D39" "3.70"
$ws.Range("E39").Value = "  +1.57%  "

# Row 40 - LidoDAOToken (D unchanged)
$ws.Range("E40").Value = "  -1.01%  "

# Row 41 - ARBITRUM
Set-PriceText "D41" "1.56"
$ws.Range("E41").Value = "  +3.75%  "

# Row 42 - BitcoinSV
Set-PriceText "D42" "98.32"
$ws.Range("E42").Value = "  +6.57%  "

# Row 43 - MultiversX
Set-PriceText "D43" "68.69"
$ws.Range("E43").Value = "  +1.26%  "

# Row 44 - FirstDigitalUSD
Set-PriceText "D44" "1.00"
$ws.Range("E44").Value = "  -0.08%  "

# Row 45 - Algorand (D unchanged)
$ws.Range("E45").Value = "  +0.87%  "

# Row 46 - Maker
Set-PriceText "D46" "1.722.01"
$ws.Range("E46").Value = "  +7.34%  "

# Row 47 - Celestia
Set-PriceText "D47" "11.88"
$ws.Range("E47").Value = "  +0.20%  "

# Row 48 - Aave
Set-PriceText "D48" "109.86"
$ws.Range("E48").Value = "  -4.61%  "

# Row 49 - ordi
Set-PriceText "D49" "77.27"
$ws.Range("E49").Value = "  -2.02%  "

# Row 50 - THORChain (D unchanged)
$ws.Range("E50").Value = "  -0.65%  "

# Row 51 - FraxShare (D unchanged)
$ws.Range("E51").Value = "  -3.31%  "
